$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 78

# Column A holds a date-formatted string ("2020-08-16") that must stay plain
# text (matching every other row in the "Fecha" column), not get
# auto-converted into a serial date. Stage it on a scratch cell far outside
# the used range (forced to Text via NumberFormat), copy only its *value*
# (so the new row doesn't inherit the scratch cell's number format) onto
# A78, then remove the scratch row entirely so it leaves no trace.
$scratch = $ws.Cells.Item(1000, 1)
$scratch.NumberFormat = "@"
$scratch.Value = "2020-08-16"
$scratch.Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163) # xlPasteValues
$scratch.EntireRow.Delete()

$ws.Cells.Item($row, 2).Value = 522162
$ws.Cells.Item($row, 3).Value = 573723
$ws.Cells.Item($row, 4).Value = 81046
$ws.Cells.Item($row, 5).Value = 56757
$ws.Cells.Item($row, 6).Value = 26.21
